# Rename the header in A1 from "param" to "name"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"

# Move the active selection to D14 (matches the recorded view state after the edit)
$ws.Range("D14").Select()
